$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48; existing rows 48:115 shift down to 49:116,
# carrying all of their original data and formatting with them (so the former
# row 115 ends up duplicated as the new row 116, matching the target diff).
$ws.Rows("48:48").Insert()

# Populate the newly inserted row 48 with the new data record.
$ws.Range("A48").Value = 11
$ws.Range("B48").Value = "Vega Monumental Concepción"
$ws.Range("C48").Value = "Bíobío"
$ws.Range("D48").Value = 44665
$ws.Range("E48").Value = 8
$ws.Range("F48").Value = 100112032
$ws.Range("G48").Value = "Zapallo italiano"
$ws.Range("H48").Value = "Sin especificar"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 200
$ws.Range("K48").Value = 11000
$ws.Range("L48").Value = 12000
$ws.Range("M48").Value = 11500
$ws.Range("N48").Value = "$/caja 50 unidades"
$ws.Range("O48").Value = "Región de Arica y Parinacota"
$ws.Range("P48").Value = 230
$ws.Range("Q48").Value = 50
$ws.Range("R48").Value = "Hortaliza"

# Make sure the date cell keeps the date number format used by the rest of
# column D (the Insert() above already carries this down from row 47, but set
# it explicitly to be safe).
$ws.Range("D48").NumberFormat = $ws.Range("D47").NumberFormat
